$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New LeetCode entry: 1920. Build Array from Permutation
$ws.Range("A5").Value = 1920
$ws.Range("B5").Value = "Build Array from Permutation"

# Header row (A1:B1) goes bold; the "#" column header (A1) is also centered
$ws.Range("A1:B1").Font.Bold = $true

# Whole "#" id column (A1:A5, including the new row) is centered
$ws.Range("A1:A5").HorizontalAlignment = -4108  # xlCenter

# Leave the active selection where the author's cursor ended up
$ws.Range("B9").Select()
